# Update the stale plate number in row 2 to the new detection result
# (KA01MG1234 -> MH20EJ0364), then leave the selection where the user's
# cursor ended up after making the edit (cell I4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "MH20EJ0364"

$ws.Range("I4").Select()
